$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6997039914131165
$ws.Range("B1").Value = 2.319869518280029
$ws.Range("C1").Value = 3.36503529548645
$ws.Range("D1").Value = 2.650514841079712
$ws.Range("E1").Value = 1.549545645713806
